$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns B through G (dates 18_12_2023, 05_01_2024, 15_01_2024,
# 21_01_2024, 28_01_2024, 07_02_2024), shifting the remaining two data
# columns (11_02_2024, 18_02_2024) left into B and C.
$ws.Range("B1:G1").EntireColumn.Delete()

$ws.Range("B:C").Select()
